$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) values. These cells store numeric-looking
# text (e.g. "245.95"), so format the cells as Text before assigning,
# then restore the default "Normal" style so no stray number format
# lingers on the cell, matching the original plain-text cell styling.
$priceUpdates = @{
    "D2" = "245.95"
    "D3" = "24.18"
    "D4" = "5.361"
    "D5" = "0.05741"
    "D6" = "6.475"
    "D7" = "3.131"
    "D8" = "0.8167"
    "D9" = "0.8713"
    "D10" = "0.1376"
    "D11" = "0.06982"
    "D12" = "0.03140"
    "D13" = "0.02932"
    "D14" = "0.09407"
    "D15" = "3.731"
    "D16" = "0.001524"
    "D17" = "0.04694"
    "D18" = "0.0005967"
    "D19" = "0.006157"
    "D20" = "0.001244"
    "D21" = "0.003908"
    "D22" = "0.00008795"
    "D23" = "3.537"
    "D24" = "2.141"
    "D25" = "0.3165"
    "D27" = "0.1327"
    "D28" = "0.0002330"
    "D40" = "0.03714"
    "D41" = "0.006406"
    "D43" = "0.002998"
    "D44" = "0.008714"
    "D45" = "0.00005269"
    "D46" = "0.00000000750"
    "D47" = "0.3898"
    "D48" = "0.002243"
    "D49" = "0.00002099"
    "D50" = "0.0001999"
}
foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}

# Updated "Volume(1h)" (column E) label text - plain text, assign directly.
$volumeUpdates = @{
    "E18" = "17OneONE"
    "E28" = "27UpBotsUBXT"
    "E44" = "43LocalTradersLCTBestin24h"
    "E48" = "47BOLOBOLOWorstin24h"
}
foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}
